$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (the sheet already has a header row 1 plus
# data rows 2.. ; this becomes the newest weekly record, pushing all the
# older records down by one row - matching the diff where every existing
# row r (>=4) became row r+1, and a brand new row 25 appeared at the end).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with this week's record.
$ws.Cells.Item(4, 1).Value = 11
$ws.Cells.Item(4, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(4, 3).Value = "Bíobío"
$ws.Cells.Item(4, 4).Value = (Get-Date -Year 2023 -Month 12 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = "Fruta"
$ws.Cells.Item(4, 7).Value = 100101
$ws.Cells.Item(4, 8).Value = "Berries"
$ws.Cells.Item(4, 9).Value = 100101004
$ws.Cells.Item(4, 10).Value = "Frambuesa"
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 220
$ws.Cells.Item(4, 14).Value = 8500
$ws.Cells.Item(4, 15).Value = 9000
$ws.Cells.Item(4, 16).Value = 8727
$ws.Cells.Item(4, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(4, 18).Value = "Región de Ñuble"
$ws.Cells.Item(4, 19).Value = 4364
$ws.Cells.Item(4, 20).Value = 2
